$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item("Content Placeholder 2")
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3)

# Original paragraph 3 text (211 chars) is made up of these runs:
#  1-100   "If you want to modify them, you must use the same reference frame and save them with the same name. "
#  101     "T"
#  102-179 "hen you must substitute the old geometries and import the PISAIIT_HAND.cmd in "
#  180-184 "adams"
#  185-211 " to generate the .bin file."
#
# Target text for that region:
#  "hen you must substitute the old geometry files and import the PISAIIT_HAND.cmd in ADAMS to generate the .bin file."
#
# Edit strictly right-to-left so earlier (still-untouched) character offsets stay valid.

# 1) " to generate the .bin file." -> "generate the .bin file."
$para.Characters(185, 27).Text = "generate the .bin file."

# 2) "in " (tail of the big run) + "adams" -> "in ADAMS to "
$para.Characters(177, 8).Text = "in ADAMS to "

# 3) "geometries and import the PISAIIT_HAND.cmd " -> "geometry files and import the PISAIIT_HAND.cmd "
$para.Characters(134, 43).Text = "geometry files and import the PISAIIT_HAND.cmd "

# 4) split that new text into two runs: "geometry files " | "and import the PISAIIT_HAND.cmd "
$para.Characters(134, 15).Text = "geometry files "
